$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input values that drive the formula recalculations
$ws.Range("B3").Value = 0.0641
$ws.Range("B4").Value = 600

# Move the active cell selection from B4 to B7
$ws.Range("B7").Select()

# Adjust the saved window x-position in the workbook view
$excel.ActiveWindow.Left = 380
